$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-16 (row 1 is header, row 4 unchanged)
$data = @(
    @("vgr_1", 1866, 1, "/vgr/pick_up_and_transport"),
    @("sm_1", 378, 2, "/sm/sort, /sm/transport"),
    @("sm_2", 309, 2, "/sm/sort, /sm/transport"),
    @("hw_1", 522, 1, "/hw/human_review"),
    @("pm_1", 204, 3, "/pm/punch_gill, /pm/punch_recesses, /pm/punch_ribbing"),
    @("mm_2", 381, 4, "/mm/deburr, /mm/drill, /mm/mill, /mm/transport_from_to"),
    @("hbw_1", 873, 4, "/hbw/get_empty_bucket, /hbw/store, /hbw/store_empty_bucket, /hbw/unload"),
    @("ov_1", 612, 2, "/ov/burn, /ov/temper"),
    @("vgr_2", 885, 1, "/vgr/pick_up_and_transport"),
    @("ov_2", 330, 1, "/ov/burn"),
    @("mm_1", 576, 4, "/mm/deburr, /mm/drill, /mm/mill, /mm/transport_from_to"),
    @("dm_2", 177, 3, "/dm/cylindrical_drill, /dm/drill, /dm/lower"),
    @("wt_1", 447, 1, "/wt/pick_up_and_transport"),
    @("wt_2", 330, 1, "/wt/pick_up_and_transport"),
    @("hbw_2", 1581, 2, "/hbw/store_empty_bucket, /hbw/unload")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
